# "Notas das listas 03 e 04"
# Adds Lista 03 / Lista 04 scores and a "Média final" (AVERAGE) column to
# the "Listas" worksheet, and makes "Listas" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Listas")

# New grades for Lista 03 (col D) and Lista 04 (col E), rows 7-26.
$grades = @{
     7 = @(100, 100)
     8 = @(0,   0)
     9 = @(90,  20)
    10 = @(80,  85)
    11 = @(0,   0)
    12 = @(90,  99)
    13 = @(0,   0)
    14 = @(0,   0)
    15 = @(80,  100)
    16 = @(90,  100)
    17 = @(0,   0)
    18 = @(0,   0)
    19 = @(70,  85)
    20 = @(90,  100)
    21 = @(75,  85)
    22 = @(85,  100)
    23 = @(95,  85)
    24 = @(0,   0)
    25 = @(0,   0)
    26 = @(90,  100)
}

# Header row: copy the header formatting from D6 (bold + centered) onto the
# new E6 / G6 header cells before writing their text.
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("E6").Value = "Lista 04"
$ws.Range("G6").Value = "Média final"

# Lista 03 / Lista 04 score data.
foreach ($r in 7..26) {
    $vals = $grades[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("E$r").Value = $vals[1]
}

# "Média final" column: AVERAGE across the four list grades, filled down
# (produces a shared formula, same as Excel's own fill-down).
$ws.Range("G7").Formula = "=AVERAGE(B7:E7)"
$ws.Range("G8:G26").Formula = "=AVERAGE(B8:E8)"

# Size column G to fit its new contents.
$ws.Columns("G").ColumnWidth = 9.67

# "Listas" becomes the selected/active sheet with the last edited cell
# selected, matching the saved view state.
$ws.Activate()
$ws.Range("L7").Select()

Write-Output "Lista 03 / Lista 04 grades and Media final column added."
